$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13.9064, 0.2606, 1.926,  1.9842, 0.3761),
    @(21.8247, 0.3299, 1.9468, 0.9789, 0.1746),
    @(6.2136,  0.1743, 1.3003, 3.8979, 0.6026),
    @(17.5509, 0.104,  0.5616, 6.3853, 0.4442),
    @(8.047599999999999, 0.1201, 5.8881, 0.6037, 0.4569)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
